$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.639.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.597.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.594.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.639.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("E21").Value = '  +5.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0508'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.626'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.843'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.50'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("E41").Value = '  +0.87%  '
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '63.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.942'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +17.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.735.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("E48").Value = '  +4.45%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.61%  '
